$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "20.611.58"
$ws.Range("E2").Value = "  +2.49%  "
$ws.Range("D3").Value = "1.469.59"
$ws.Range("E3").Value = "  +2.77%  "
$ws.Range("E4").Value = "  +0.98%  "
$ws.Range("D5").Value = "'0.9588"
$ws.Range("E5").Value = "  -4.18%  "
$ws.Range("D6").Value = "'282.01"
$ws.Range("E6").Value = "  +2.55%  "
$ws.Range("D7").Value = "'0.3731"
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("D8").Value = "'0.3203"
$ws.Range("E8").Value = "  +3.60%  "
$ws.Range("D9").Value = "'41.92"
$ws.Range("E9").Value = "  +4.38%  "
$ws.Range("D10").Value = "'1.067"
$ws.Range("E10").Value = "  +5.29%  "
$ws.Range("D11").Value = "'0.06723"
$ws.Range("E11").Value = "  +1.85%  "
$ws.Range("D12").Value = "'1.003"
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("D13").Value = "'5.651"
$ws.Range("E13").Value = "  +4.48%  "
$ws.Range("D14").Value = "'18.42"
$ws.Range("E14").Value = "  +6.73%  "
$ws.Range("D15").Value = "'6.297"
$ws.Range("E15").Value = "  +1.79%  "
$ws.Range("D16").Value = "1.476.74"
$ws.Range("E16").Value = "  +3.33%  "
$ws.Range("E17").Value = "  +2.88%  "
$ws.Range("D18").Value = "'0.05792"
$ws.Range("E18").Value = "  -0.45%  "
$ws.Range("D19").Value = "'73.00"
$ws.Range("E19").Value = "  -3.59%  "
$ws.Range("D20").Value = "'0.9598"
$ws.Range("E20").Value = "  -4.06%  "
$ws.Range("D21").Value = "'5.724"
$ws.Range("E21").Value = "  +0.48%  "
$ws.Range("D22").Value = "'14.91"
$ws.Range("E22").Value = "  +2.55%  "
$ws.Range("D23").Value = "'11.26"
$ws.Range("E23").Value = "  +1.50%  "
$ws.Range("D24").Value = "'2.302"
$ws.Range("E24").Value = "  -1.45%  "
$ws.Range("D25").Value = "20.719.21"
$ws.Range("E25").Value = "  +2.95%  "
$ws.Range("D26").Value = "'2.333"
$ws.Range("E26").Value = "  +1.64%  "
$ws.Range("D27").Value = "'137.94"
$ws.Range("E27").Value = "  -0.36%  "
$ws.Range("D28").Value = "'17.66"
$ws.Range("E28").Value = "  +4.32%  "
$ws.Range("D29").Value = "1.637.81"
$ws.Range("E29").Value = "  +2.91%  "
$ws.Range("D30").Value = "'114.01"
$ws.Range("E30").Value = "  +4.08%  "
$ws.Range("D31").Value = "'3.973"
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("D32").Value = "'5.395"
$ws.Range("E32").Value = "  -0.67%  "
$ws.Range("D33").Value = "'0.8463"
$ws.Range("E33").Value = "  -7.49%  "
$ws.Range("D34").Value = "'1.647"
$ws.Range("E34").Value = "  +26.72%  "
$ws.Range("D35").Value = "'0.07872"
$ws.Range("E35").Value = "  +1.24%  "
$ws.Range("D36").Value = "'0.06118"
$ws.Range("E36").Value = "  +7.21%  "
$ws.Range("D37").Value = "'4.965"
$ws.Range("E37").Value = "  +4.21%  "
$ws.Range("D38").Value = "'10.78"
$ws.Range("E38").Value = "  -6.06%  "
$ws.Range("D39").Value = "'0.02085"
$ws.Range("E39").Value = "  +2.77%  "
$ws.Range("D40").Value = "'1.134"
$ws.Range("E40").Value = "  +1.00%  "
$ws.Range("D41").Value = "'0.9694"
$ws.Range("E41").Value = "  -3.07%  "
$ws.Range("D42").Value = "'0.1911"
$ws.Range("E42").Value = "  -0.77%  "
$ws.Range("D43").Value = "'7.525"
$ws.Range("E43").Value = "  -10.92%  "
$ws.Range("D44").Value = "'0.5468"
$ws.Range("E44").Value = "  +2.15%  "
$ws.Range("D45").Value = "'12.59"
$ws.Range("E45").Value = "  +2.84%  "
$ws.Range("D46").Value = "'3.598"
$ws.Range("E46").Value = "  +1.32%  "
$ws.Range("D47").Value = "'121.76"
$ws.Range("E47").Value = "  +10.84%  "
$ws.Range("D48").Value = "'0.5395"
$ws.Range("E48").Value = "  +4.72%  "
$ws.Range("E49").Value = "  +3.43%  "
$ws.Range("D50").Value = "'0.06457"
$ws.Range("E50").Value = "  +4.02%  "
$ws.Range("D51").Value = "'1.058"
$ws.Range("E51").Value = "  +0.40%  "
